$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new timesheet entry for row 6
$ws.Range("A6").Value = 42627
$ws.Range("A6").NumberFormat = "d-mmm"

$ws.Range("B6").Value = 0.65277777777777779
$ws.Range("B6").NumberFormat = "h:mm"

$ws.Range("C6").Value = 0.68055555555555547
$ws.Range("C6").NumberFormat = "h:mm"

$ws.Range("D6").Value = "Implementação do cadastrar cliente"

# Move the active selection to D7, matching the saved view state
$ws.Range("D7").Select()
